# Insert a new weekly price record for "Piña" (Macroferia Regional de Talca)
# as row 164, pushing the existing rows 164-217 down to 165-218.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 164..217 down by one, carrying formatting (incl. the date
# style on column D) from the row being pushed down.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new observation.
$ws.Cells.Item(164, 1).Value  = 5
$ws.Cells.Item(164, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(164, 3).Value  = "Maule"
$ws.Cells.Item(164, 4).Value  = 44627
$ws.Cells.Item(164, 5).Value  = 7
$ws.Cells.Item(164, 6).Value  = "Fruta"
$ws.Cells.Item(164, 7).Value  = 100108
$ws.Cells.Item(164, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(164, 9).Value  = 100108005
$ws.Cells.Item(164, 10).Value = "Piña"
$ws.Cells.Item(164, 11).Value = "Caramelo"
$ws.Cells.Item(164, 12).Value = "Segunda"
$ws.Cells.Item(164, 13).Value = 180
$ws.Cells.Item(164, 14).Value = 19000
$ws.Cells.Item(164, 15).Value = 19000
$ws.Cells.Item(164, 16).Value = 19000
$ws.Cells.Item(164, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(164, 18).Value = "Ecuador"
$ws.Cells.Item(164, 19).Value = 1357
$ws.Cells.Item(164, 20).Value = 14
